# Pacific Cross Master plan: switch the quoted currency from USD to VND
# and reset the conversion rate to 1 (VND is now the base/local currency,
# so no conversion multiplier is needed). Finally leave the selection
# where the user's cursor would land after editing the currency cell
# (one row below it, at E3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "VND"
$ws.Range("C2").Value = 1

$ws.Range("E3").Select()
